$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells D2:E51 store plain text values (e.g. "26.342.08", "  +0.20%  ").
# Force text format before assigning so Excel does not reinterpret
# strings like "218.15" or "0.5423" as numbers.
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "26.368.86"
$ws.Cells.Item(2, 5).Value = "  +0.39%  "

$ws.Cells.Item(3, 4).Value = "1.696.59"
$ws.Cells.Item(3, 5).Value = "  +1.06%  "

$ws.Cells.Item(4, 5).Value = "  +0.07%  "

$ws.Cells.Item(5, 4).Value = "218.15"
$ws.Cells.Item(5, 5).Value = "  -0.07%  "

$ws.Cells.Item(6, 4).Value = "0.5423"
$ws.Cells.Item(6, 5).Value = "  +2.87%  "

$ws.Cells.Item(7, 5).Value = "  +0.01%  "

$ws.Cells.Item(8, 5).Value = "  +1.17%  "

$ws.Cells.Item(9, 4).Value = "0.06466"
$ws.Cells.Item(9, 5).Value = "  -0.08%  "

$ws.Cells.Item(10, 5).Value = "  -1.25%  "

$ws.Cells.Item(11, 4).Value = "0.07667"
$ws.Cells.Item(11, 5).Value = "  +1.92%  "

$ws.Cells.Item(12, 4).Value = "1.689.12"
$ws.Cells.Item(12, 5).Value = "  +0.53%  "

$ws.Cells.Item(13, 4).Value = "4.546"
$ws.Cells.Item(13, 5).Value = "  +0.39%  "

$ws.Cells.Item(14, 4).Value = "0.5830"
$ws.Cells.Item(14, 5).Value = "  +0.55%  "

$ws.Cells.Item(15, 4).Value = "0.000008426"
$ws.Cells.Item(15, 5).Value = "  -0.87%  "

$ws.Cells.Item(16, 5).Value = "  +3.60%  "

$ws.Cells.Item(17, 4).Value = "26.407.62"
$ws.Cells.Item(17, 5).Value = "  +0.34%  "

$ws.Cells.Item(18, 4).Value = "4.926"
$ws.Cells.Item(18, 5).Value = "  +0.16%  "

$ws.Cells.Item(19, 5).Value = "  +0.06%  "

$ws.Cells.Item(20, 5).Value = "  +0.39%  "

$ws.Cells.Item(21, 4).Value = "191.33"
$ws.Cells.Item(21, 5).Value = "  +0.57%  "

$ws.Cells.Item(22, 4).Value = "6.299"
$ws.Cells.Item(22, 5).Value = "  +1.72%  "

$ws.Cells.Item(23, 5).Value = "  +0.14%  "

$ws.Cells.Item(24, 4).Value = "149.02"
$ws.Cells.Item(24, 5).Value = "  +2.75%  "

$ws.Cells.Item(25, 4).Value = "0.1292"
$ws.Cells.Item(25, 5).Value = "  +3.97%  "

$ws.Cells.Item(26, 4).Value = "7.877"
$ws.Cells.Item(26, 5).Value = "  +0.87%  "

$ws.Cells.Item(27, 4).Value = "15.90"
$ws.Cells.Item(27, 5).Value = "  +0.57%  "

$ws.Cells.Item(28, 4).Value = "0.06373"
$ws.Cells.Item(28, 5).Value = "  -2.96%  "

$ws.Cells.Item(29, 4).Value = "1.389"
$ws.Cells.Item(29, 5).Value = "  +2.48%  "

$ws.Cells.Item(30, 5).Value = "  -0.10%  "

$ws.Cells.Item(31, 4).Value = "3.619"
$ws.Cells.Item(31, 5).Value = "  +0.90%  "

$ws.Cells.Item(32, 4).Value = "3.596"
$ws.Cells.Item(32, 5).Value = "  -0.01%  "

$ws.Cells.Item(33, 4).Value = "1.694"
$ws.Cells.Item(33, 5).Value = "  +2.11%  "

$ws.Cells.Item(34, 4).Value = "1.036"
$ws.Cells.Item(34, 5).Value = "  +0.32%  "

$ws.Cells.Item(35, 4).Value = "0.6211"
$ws.Cells.Item(35, 5).Value = "  -0.22%  "

$ws.Cells.Item(36, 4).Value = "2.416"
$ws.Cells.Item(36, 5).Value = "  +0.50%  "

$ws.Cells.Item(37, 4).Value = "2.752"
$ws.Cells.Item(37, 5).Value = "  +0.30%  "

$ws.Cells.Item(38, 4).Value = "0.01657"
$ws.Cells.Item(38, 5).Value = "  +2.14%  "

$ws.Cells.Item(39, 4).Value = "1.117.84"
$ws.Cells.Item(39, 5).Value = "  +0.47%  "

$ws.Cells.Item(40, 4).Value = "6.105"
$ws.Cells.Item(40, 5).Value = "  -5.31%  "

$ws.Cells.Item(41, 4).Value = "0.8848"
$ws.Cells.Item(41, 5).Value = "  +0.89%  "

$ws.Cells.Item(42, 5).Value = "  -0.01%  "

$ws.Cells.Item(43, 4).Value = "101.11"
$ws.Cells.Item(43, 5).Value = "  +0.33%  "

$ws.Cells.Item(44, 4).Value = "1.847.26"
$ws.Cells.Item(44, 5).Value = "  +1.01%  "

$ws.Cells.Item(45, 4).Value = "0.00000000111"
$ws.Cells.Item(45, 5).Value = "  -1.85%  "

$ws.Cells.Item(46, 4).Value = "57.81"
$ws.Cells.Item(46, 5).Value = "  +1.50%  "

$ws.Cells.Item(47, 4).Value = "8.202"
$ws.Cells.Item(47, 5).Value = "  +0.09%  "

$ws.Cells.Item(48, 4).Value = "1.002"
$ws.Cells.Item(48, 5).Value = "  -0.40%  "

$ws.Cells.Item(49, 5).Value = "  +0.30%  "

$ws.Cells.Item(50, 4).Value = "6.112"
$ws.Cells.Item(50, 5).Value = "  +0.23%  "

$ws.Cells.Item(51, 5).Value = "  +0.17%  "

# Restore original (default) cell formatting now that the text values are set,
# so no stray number-format styling is left behind on these cells.
$rng.ClearFormats()
